$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range('A2').Value = 'Última actualización: 11:03:46'
$ws.Range('A3').Value = 'Total filas: 163'

$ws.Cells.Item(44, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(46, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(86, 1).Value = '07:31:43'
$ws.Cells.Item(86, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(86, 4).Value = 99
$ws.Cells.Item(87, 1).Value = '08:42:31'
$ws.Cells.Item(87, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(87, 4).Value = 28
$ws.Cells.Item(96, 1).Value = '07:31:43'
$ws.Cells.Item(96, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(96, 4).Value = 112
$ws.Cells.Item(97, 1).Value = '07:57:27'
$ws.Cells.Item(97, 3).Value = '17_ROMERO'
$ws.Cells.Item(97, 4).Value = 86
$ws.Cells.Item(132, 1).Value = '11:03:46'
$ws.Cells.Item(132, 2).Value = '11:03'
$ws.Cells.Item(132, 3).Value = '215C_EL PATO'
$ws.Cells.Item(132, 4).Value = 0
$ws.Cells.Item(132, 5).Value = 'LP1912'
$ws.Cells.Item(133, 1).Value = '10:30:21'
$ws.Cells.Item(133, 2).Value = '11:03'
$ws.Cells.Item(133, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(133, 4).Value = 33
$ws.Cells.Item(133, 5).Value = 'LP1912'
$ws.Cells.Item(134, 1).Value = '11:03:46'
$ws.Cells.Item(134, 2).Value = '11:04'
$ws.Cells.Item(134, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(134, 4).Value = 1
$ws.Cells.Item(134, 5).Value = 'LP1912'
$ws.Cells.Item(135, 1).Value = '11:03:46'
$ws.Cells.Item(135, 2).Value = '11:04'
$ws.Cells.Item(135, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(135, 4).Value = 1
$ws.Cells.Item(135, 5).Value = 'LP1912'
$ws.Cells.Item(136, 1).Value = '10:30:21'
$ws.Cells.Item(136, 2).Value = '11:06'
$ws.Cells.Item(136, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(136, 4).Value = 36
$ws.Cells.Item(136, 5).Value = 'LP1912'
$ws.Cells.Item(137, 1).Value = '10:30:21'
$ws.Cells.Item(137, 2).Value = '11:11'
$ws.Cells.Item(137, 3).Value = '10_OLMOS'
$ws.Cells.Item(137, 4).Value = 41
$ws.Cells.Item(137, 5).Value = 'LP1912'
$ws.Cells.Item(138, 1).Value = '11:03:46'
$ws.Cells.Item(138, 2).Value = '11:12'
$ws.Cells.Item(138, 3).Value = '15_ABASTO'
$ws.Cells.Item(138, 4).Value = 9
$ws.Cells.Item(138, 5).Value = 'LP1912'
$ws.Cells.Item(139, 1).Value = '09:31:25'
$ws.Cells.Item(139, 2).Value = '11:17'
$ws.Cells.Item(139, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(139, 4).Value = 106
$ws.Cells.Item(139, 5).Value = 'LP1912'
$ws.Cells.Item(140, 1).Value = '09:31:25'
$ws.Cells.Item(140, 2).Value = '11:19'
$ws.Cells.Item(140, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(140, 4).Value = 108
$ws.Cells.Item(140, 5).Value = 'LP1912'
$ws.Cells.Item(141, 1).Value = '09:31:25'
$ws.Cells.Item(141, 2).Value = '11:21'
$ws.Cells.Item(141, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(141, 4).Value = 110
$ws.Cells.Item(141, 5).Value = 'LP1912'
$ws.Cells.Item(142, 1).Value = '09:31:25'
$ws.Cells.Item(142, 2).Value = '11:26'
$ws.Cells.Item(142, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(142, 4).Value = 115
$ws.Cells.Item(142, 5).Value = 'LP1912'
$ws.Cells.Item(143, 1).Value = '09:31:25'
$ws.Cells.Item(143, 2).Value = '11:27'
$ws.Cells.Item(143, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(143, 4).Value = 116
$ws.Cells.Item(143, 5).Value = 'LP1912'
$ws.Cells.Item(144, 1).Value = '10:30:21'
$ws.Cells.Item(144, 2).Value = '11:32'
$ws.Cells.Item(144, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(144, 4).Value = 62
$ws.Cells.Item(144, 5).Value = 'LP1912'
$ws.Cells.Item(145, 1).Value = '11:03:46'
$ws.Cells.Item(145, 2).Value = '11:34'
$ws.Cells.Item(145, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(145, 4).Value = 31
$ws.Cells.Item(145, 5).Value = 'LP1912'
$ws.Cells.Item(146, 1).Value = '10:30:21'
$ws.Cells.Item(146, 2).Value = '11:35'
$ws.Cells.Item(146, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(146, 4).Value = 65
$ws.Cells.Item(146, 5).Value = 'LP1912'
$ws.Cells.Item(147, 1).Value = '10:30:21'
$ws.Cells.Item(147, 2).Value = '11:39'
$ws.Cells.Item(147, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(147, 4).Value = 69
$ws.Cells.Item(147, 5).Value = 'LP1912'
$ws.Cells.Item(148, 1).Value = '11:03:46'
$ws.Cells.Item(148, 2).Value = '11:41'
$ws.Cells.Item(148, 3).Value = '17_ROMERO'
$ws.Cells.Item(148, 4).Value = 38
$ws.Cells.Item(148, 5).Value = 'LP1912'
$ws.Cells.Item(149, 1).Value = '10:30:21'
$ws.Cells.Item(149, 2).Value = '11:42'
$ws.Cells.Item(149, 3).Value = '17_ROMERO'
$ws.Cells.Item(149, 4).Value = 72
$ws.Cells.Item(149, 5).Value = 'LP1912'
$ws.Cells.Item(150, 1).Value = '11:03:46'
$ws.Cells.Item(150, 2).Value = '11:43'
$ws.Cells.Item(150, 3).Value = '10_OLMOS'
$ws.Cells.Item(150, 4).Value = 40
$ws.Cells.Item(150, 5).Value = 'LP1912'
$ws.Cells.Item(151, 1).Value = '10:30:21'
$ws.Cells.Item(151, 2).Value = '11:48'
$ws.Cells.Item(151, 3).Value = '10_OLMOS'
$ws.Cells.Item(151, 4).Value = 78
$ws.Cells.Item(151, 5).Value = 'LP1912'
$ws.Cells.Item(152, 1).Value = '10:30:21'
$ws.Cells.Item(152, 2).Value = '11:51'
$ws.Cells.Item(152, 3).Value = '215B_EL PATO'
$ws.Cells.Item(152, 4).Value = 81
$ws.Cells.Item(152, 5).Value = 'LP1912'
$ws.Cells.Item(153, 1).Value = '11:03:46'
$ws.Cells.Item(153, 2).Value = '11:52'
$ws.Cells.Item(153, 3).Value = '15_ABASTO'
$ws.Cells.Item(153, 4).Value = 49
$ws.Cells.Item(153, 5).Value = 'LP1912'
$ws.Cells.Item(154, 1).Value = '10:30:21'
$ws.Cells.Item(154, 2).Value = '11:54'
$ws.Cells.Item(154, 3).Value = '15_ABASTO'
$ws.Cells.Item(154, 4).Value = 84
$ws.Cells.Item(154, 5).Value = 'LP1912'
$ws.Cells.Item(155, 1).Value = '10:30:21'
$ws.Cells.Item(155, 2).Value = '11:59'
$ws.Cells.Item(155, 3).Value = '225_GOMEZ'
$ws.Cells.Item(155, 4).Value = 89
$ws.Cells.Item(155, 5).Value = 'LP1912'
$ws.Cells.Item(156, 1).Value = '10:30:21'
$ws.Cells.Item(156, 2).Value = '12:02'
$ws.Cells.Item(156, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(156, 4).Value = 92
$ws.Cells.Item(156, 5).Value = 'LP1912'
$ws.Cells.Item(157, 1).Value = '10:30:21'
$ws.Cells.Item(157, 2).Value = '12:06'
$ws.Cells.Item(157, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(157, 4).Value = 96
$ws.Cells.Item(157, 5).Value = 'LP1912'
$ws.Cells.Item(158, 1).Value = '11:03:46'
$ws.Cells.Item(158, 2).Value = '12:06'
$ws.Cells.Item(158, 3).Value = '14_ABASTO'
$ws.Cells.Item(158, 4).Value = 63
$ws.Cells.Item(158, 5).Value = 'LP1912'
$ws.Cells.Item(159, 1).Value = '11:03:46'
$ws.Cells.Item(159, 2).Value = '12:09'
$ws.Cells.Item(159, 3).Value = '10_OLMOS'
$ws.Cells.Item(159, 4).Value = 66
$ws.Cells.Item(159, 5).Value = 'LP1912'
$ws.Cells.Item(160, 1).Value = '10:30:21'
$ws.Cells.Item(160, 2).Value = '12:14'
$ws.Cells.Item(160, 3).Value = '17_ROMERO'
$ws.Cells.Item(160, 4).Value = 104
$ws.Cells.Item(160, 5).Value = 'LP1912'
$ws.Cells.Item(161, 1).Value = '10:30:21'
$ws.Cells.Item(161, 2).Value = '12:17'
$ws.Cells.Item(161, 3).Value = '14_ABASTO'
$ws.Cells.Item(161, 4).Value = 107
$ws.Cells.Item(161, 5).Value = 'LP1912'
$ws.Cells.Item(162, 1).Value = '10:30:21'
$ws.Cells.Item(162, 2).Value = '12:20'
$ws.Cells.Item(162, 3).Value = '215A_EL PATO'
$ws.Cells.Item(162, 4).Value = 110
$ws.Cells.Item(162, 5).Value = 'LP1912'
$ws.Cells.Item(163, 1).Value = '11:03:46'
$ws.Cells.Item(163, 2).Value = '12:20'
$ws.Cells.Item(163, 3).Value = '14_ABASTO'
$ws.Cells.Item(163, 4).Value = 77
$ws.Cells.Item(163, 5).Value = 'LP1912'
$ws.Cells.Item(164, 1).Value = '10:30:21'
$ws.Cells.Item(164, 2).Value = '12:21'
$ws.Cells.Item(164, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(164, 4).Value = 111
$ws.Cells.Item(164, 5).Value = 'LP1912'
$ws.Cells.Item(165, 1).Value = '11:03:46'
$ws.Cells.Item(165, 2).Value = '12:31'
$ws.Cells.Item(165, 3).Value = '17_ROMERO'
$ws.Cells.Item(165, 4).Value = 88
$ws.Cells.Item(165, 5).Value = 'LP1912'
$ws.Cells.Item(166, 1).Value = '11:03:46'
$ws.Cells.Item(166, 2).Value = '12:36'
$ws.Cells.Item(166, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(166, 4).Value = 93
$ws.Cells.Item(166, 5).Value = 'LP1912'
$ws.Cells.Item(167, 1).Value = '11:03:46'
$ws.Cells.Item(167, 2).Value = '12:38'
$ws.Cells.Item(167, 3).Value = '17_179 Y 38'
$ws.Cells.Item(167, 4).Value = 95
$ws.Cells.Item(167, 5).Value = 'LP1912'
$ws.Cells.Item(168, 1).Value = '11:03:46'
$ws.Cells.Item(168, 2).Value = '12:48'
$ws.Cells.Item(168, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(168, 4).Value = 105
$ws.Cells.Item(168, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range('A2').Value = 'Última actualización: 11:03:46'
$ws.Range('A3').Value = 'Total filas: 18'

$ws.Cells.Item(21, 1).Value = '11:03:46'
$ws.Cells.Item(21, 2).Value = '11:03'
$ws.Cells.Item(21, 3).Value = '215C_EL PATO'
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 'LP1912'
$ws.Cells.Item(22, 1).Value = '10:30:21'
$ws.Cells.Item(22, 2).Value = '11:51'
$ws.Cells.Item(22, 3).Value = '215B_EL PATO'
$ws.Cells.Item(22, 4).Value = 81
$ws.Cells.Item(22, 5).Value = 'LP1912'
$ws.Cells.Item(23, 1).Value = '10:30:21'
$ws.Cells.Item(23, 2).Value = '12:20'
$ws.Cells.Item(23, 3).Value = '215A_EL PATO'
$ws.Cells.Item(23, 4).Value = 110
$ws.Cells.Item(23, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range('A2').Value = 'Última actualización: 11:03:46'
$ws.Range('A3').Value = 'Total filas: 27'

$ws.Cells.Item(29, 1).Value = '11:03:46'
$ws.Cells.Item(29, 2).Value = '11:13'
$ws.Cells.Item(29, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(29, 4).Value = 10
$ws.Cells.Item(29, 5).Value = 'L6203'
$ws.Cells.Item(30, 1).Value = '09:31:25'
$ws.Cells.Item(30, 2).Value = '11:14'
$ws.Cells.Item(30, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(30, 4).Value = 103
$ws.Cells.Item(30, 5).Value = 'L6203'
$ws.Cells.Item(31, 1).Value = '10:30:21'
$ws.Cells.Item(31, 2).Value = '12:04'
$ws.Cells.Item(31, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(31, 4).Value = 94
$ws.Cells.Item(31, 5).Value = 'L6173'
$ws.Cells.Item(32, 1).Value = '11:03:46'
$ws.Cells.Item(32, 2).Value = '12:53'
$ws.Cells.Item(32, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(32, 4).Value = 110
$ws.Cells.Item(32, 5).Value = 'L6203'
